$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Week 1" -> "Week 2" for column A data rows (2-27)
$ws.Range("A2:A27").Value = "Week 2"

# Set column widths for D and G (observed from target file)
$ws.Columns("D").ColumnWidth = 16.5
$ws.Columns("G").ColumnWidth = 21.67

# Reset the view: select H11 (clears the prior top-left scroll position)
$ws.Range("H11").Select()
